# Master_Fuel_Sector_List.xlsx update:
# Add new sector "7BC_Indirect-N2O-non-agricultural-N" to the "Sectors" sheet,
# inserted right after "7A_Fossil-fuel-fires" (i.e. as the new row 58, pushing
# "11A_Volcanoes", "11B_Forest-fires" and "11C_Other-natural" down by one row).
# It uses the same activity/units/type as the other "NC"-type sectors:
#   activity = "pop", units = 1000, type = "NC"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# Insert a new blank row above the current row 58 ("11A_Volcanoes"),
# shifting the existing rows 58-60 down to 59-61.
$ws.Rows.Item(58).Insert()

# Fill in the new sector row.
$ws.Range("A58").Value = "7BC_Indirect-N2O-non-agricultural-N"
$ws.Range("B58").Value = "pop"
$ws.Range("C58").Value = 1000
$ws.Range("D58").Value = "NC"

# Match formatting of the surrounding data rows (copy format down from row 57).
$ws.Range("A57:D57").Copy()
$ws.Range("A58:D58").PasteSpecial(-4122)

# Cosmetic cleanup that also happened in the authored edit: row 23's A cell
# ("1A3di_Oil_Tanker_Loading") had been using a stray one-off font style;
# restore it to the same style used by the rest of column A.
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)

$excel.CutCopyMode = 0
